$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.872.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +6.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.253.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.41%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'579.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.58%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'182.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +9.01%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  -1.58%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.250.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.80%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +6.19%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.817.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.18%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.37%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'28.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.82%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'67.860.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +6.36%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000169"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.248.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.05%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.90%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.91%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'375.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +7.34%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +6.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.15%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'71.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.511"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.88%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.63%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +3.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +5.68%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.70%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'22.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +7.70%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'USDe"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'164.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.38%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +6.32%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +11.39%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +14.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'26.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'dogwifhat"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'2.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +8.29%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'365.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +15.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'4.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +7.44%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.709.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.05%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'25.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +10.38%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'40.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.94%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0679"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0280"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.41%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +8.52%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.42%  "
$ws.Range("E51").Style = "Normal"
